$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This reproduces a reordering of rows 31-35 (all other rows are untouched).
# The new row order (by original row) is: 34, 33, 31, 35, 32
# i.e. new row31 <- old row34, new row32 <- old row33, new row33 <- old row31,
#      new row34 <- old row35, new row35 <- old row32.
# Columns A,B,E,F,G,H,I,J,K,M,Q,R carry the per-record data that moves with the row;
# all other columns (C,D,L,N,P,S,T,U,V,W,Y,Z,AA,AB,AD,AE,AF,AG,AT,AW,AX,AY) are identical
# across these five rows, so they do not need to change.

# Row 31 (was row 34): Grönpyrola / Pyrola chlorantha, Antal=2, m², fullt utvecklade blad, M present=No
$ws.Range("A31").Value = 111565033
$ws.Range("B31").Value = 103288
$ws.Range("E31").Value = 221144
$ws.Range("F31").Value = "Grönpyrola"
$ws.Range("G31").Value = "Pyrola chlorantha"
$ws.Range("H31").Value = "Sw."
$ws.Range("I31").Value = "'2"
$ws.Range("I31").Style = "Normal"
$ws.Range("J31").Value = "m²"
$ws.Range("K31").Value = "fullt utvecklade blad"
$ws.Range("M31").Value = ""
$ws.Range("Q31").Value = 561151.5115810917
$ws.Range("R31").Value = 6622728.260846013

# Row 32 (was row 33): Grönpyrola / Pyrola chlorantha, Antal=2, m², fullt utvecklade blad, M present=No
$ws.Range("A32").Value = 111565024
$ws.Range("B32").Value = 103288
$ws.Range("E32").Value = 221144
$ws.Range("F32").Value = "Grönpyrola"
$ws.Range("G32").Value = "Pyrola chlorantha"
$ws.Range("H32").Value = "Sw."
$ws.Range("I32").Value = "'2"
$ws.Range("I32").Style = "Normal"
$ws.Range("J32").Value = "m²"
$ws.Range("K32").Value = "fullt utvecklade blad"
$ws.Range("Q32").Value = 561149.6074341368
$ws.Range("R32").Value = 6622721.170183762

# Row 33 (was row 31): Åkergroda / Rana arvalis, Antal=1, (no unit), årsunge, M present=Yes
$ws.Range("A33").Value = 111565017
$ws.Range("B33").Value = 57578
$ws.Range("E33").Value = 208250
$ws.Range("F33").Value = "Åkergroda"
$ws.Range("G33").Value = "Rana arvalis"
$ws.Range("H33").Value = "Nilsson, 1842"
$ws.Range("I33").Value = "'1"
$ws.Range("I33").Style = "Normal"
$ws.Range("J33").Value = "'"
$ws.Range("J33").Style = "Normal"
$ws.Range("K33").Value = "årsunge"
$ws.Range("M33").Value = "'"
$ws.Range("M33").Style = "Normal"
$ws.Range("Q33").Value = 561130.0283522989
$ws.Range("R33").Value = 6622683.03052416

# Row 34 (was row 35): Åkergroda / Rana arvalis, Antal=1, ex., årsunge, M present=Yes
$ws.Range("A34").Value = 111564905
$ws.Range("B34").Value = 57578
$ws.Range("E34").Value = 208250
$ws.Range("F34").Value = "Åkergroda"
$ws.Range("G34").Value = "Rana arvalis"
$ws.Range("H34").Value = "Nilsson, 1842"
$ws.Range("I34").Value = "'1"
$ws.Range("I34").Style = "Normal"
$ws.Range("J34").Value = "ex."
$ws.Range("K34").Value = "årsunge"
$ws.Range("M34").Value = "'"
$ws.Range("M34").Style = "Normal"
$ws.Range("Q34").Value = 561262.631747936
$ws.Range("R34").Value = 6622544.013810508

# Row 35 (was row 32): Grönpyrola / Pyrola chlorantha, Antal=3, m², fullt utvecklade blad, M present=No
$ws.Range("A35").Value = 111564885
$ws.Range("B35").Value = 103288
$ws.Range("E35").Value = 221144
$ws.Range("F35").Value = "Grönpyrola"
$ws.Range("G35").Value = "Pyrola chlorantha"
$ws.Range("H35").Value = "Sw."
$ws.Range("I35").Value = "'3"
$ws.Range("I35").Style = "Normal"
$ws.Range("J35").Value = "m²"
$ws.Range("K35").Value = "fullt utvecklade blad"
$ws.Range("M35").Value = ""
$ws.Range("Q35").Value = 561241.7941118333
$ws.Range("R35").Value = 6622674.779475109
